# Weekly data refresh: insert the newest week's record for
# "Vega Modelo de Temuco - Apio" ahead of the existing row 306,
# shifting all subsequent rows down by one.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a brand-new blank row at position 306 (pushes old 306..385 -> 307..386)
$ws.Rows.Item(306).Insert()

# Populate the newly inserted row with this week's reading
$ws.Cells.Item(306, 1).Value = 10
$ws.Cells.Item(306, 2).Value = "Vega Modelo de Temuco"
$ws.Cells.Item(306, 3).Value = "La Araucanía"
$ws.Cells.Item(306, 4).Value = 44855
$ws.Cells.Item(306, 5).Value = 9
$ws.Cells.Item(306, 6).Value = 100112017
$ws.Cells.Item(306, 7).Value = "Apio"
$ws.Cells.Item(306, 8).Value = "Americana (o)"
$ws.Cells.Item(306, 9).Value = "Primera"
$ws.Cells.Item(306, 10).Value = 80
$ws.Cells.Item(306, 11).Value = 9000
$ws.Cells.Item(306, 12).Value = 10000
$ws.Cells.Item(306, 13).Value = 9375
$ws.Cells.Item(306, 14).Value = "$/docena de matas"
$ws.Cells.Item(306, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(306, 16).Value = 1562
$ws.Cells.Item(306, 17).Value = 6
$ws.Cells.Item(306, 18).Value = "Hortaliza"
